# "updated retracted in fsm"
# Insert the new FSM boolean globals ("fully_retracted", "fully_extended",
# "retracted") into the Global sheet, and the new Motor.h angle macros
# ("RET_ANG", "EXT_ANG") into the MACROS sheet.

$wb = $excel.ActiveWorkbook

# --- Global sheet: insert 3 rows at row 3 (pushes everything below down by 3) ---
$wsGlobal = $wb.Worksheets.Item("Global")
$wsGlobal.Rows("3:5").Insert()

$wsGlobal.Range("A3").Value = "FSM_State.h"
$wsGlobal.Range("B3").Value = "boolean"
$wsGlobal.Range("C3").Value = "fully_retracted"
$wsGlobal.Range("D3").Value = "global variable for keeping track of if the leg is fully retracted"

$wsGlobal.Range("A4").Value = "FSM_State.h"
$wsGlobal.Range("B4").Value = "boolean"
$wsGlobal.Range("C4").Value = "fully_extended"
$wsGlobal.Range("D4").Value = "global variable for keeping track of if the leg is fully extended"

$wsGlobal.Range("A5").Value = "FSM_State.h"
$wsGlobal.Range("B5").Value = "boolean"
$wsGlobal.Range("C5").Value = "retracted"
$wsGlobal.Range("D5").Value = "global variable for keeping track of whether the leg has been retracted recently (reset to false when leg changes out of RETRACTED State)"

# --- MACROS sheet: insert 2 rows at row 34 (pushes Thresholds.h section down by 2) ---
$wsMacros = $wb.Worksheets.Item("MACROS")
$wsMacros.Rows("34:35").Insert()

$wsMacros.Range("A34").Value = "Motor.h"
$wsMacros.Range("A35").Value = "Motor.h"
$wsMacros.Range("B34").Value = "RET_ANG"
$wsMacros.Range("B35").Value = "EXT_ANG"
$wsMacros.Range("C34").Value = "angle of fully retracted leg"
$wsMacros.Range("C35").Value = "angle of fully extended leg"

# --- Restore the selections recorded in the saved workbook views ---
$wsGlobal.Activate() | Out-Null
$wsGlobal.Range("D5").Select() | Out-Null

$wsMacros.Activate() | Out-Null
$wsMacros.Range("C35").Select() | Out-Null
